$d = $word.ActiveDocument

$replacements = @(
    @("526×5=2630", "695×8=5560"),
    @("676×6=4056", "270×7=1890"),
    @("880×6=5280", "538×6=3228"),
    @("311×8=2488", "931×8=7448"),
    @("930×5=4650", "901×5=4505"),
    @("975×5=4875", "873×4=3492"),
    @("295×9=2655", "673×2=1346"),
    @("193×2=386",  "454×3=1362"),
    @("555×7=3885", "630×6=3780"),
    @("593×5=2965", "960×6=5760"),
    @("524×9=4716", "652×5=3260"),
    @("334×6=2004", "844×6=5064"),
    @("263×2=526",  "630×5=3150"),
    @("356×4=1424", "710×8=5680"),
    @("208×7=1456", "910×5=4550"),
    @("170×3=510",  "155×9=1395"),
    @("400×2=800",  "131×7=917"),
    @("959×7=6713", "602×2=1204"),
    @("991×6=5946", "835×7=5845"),
    @("998×7=6986", "392×2=784"),
    @("114×8=912",  "187×8=1496"),
    @("760×5=3800", "268×5=1340"),
    @("952×5=4760", "627×7=4389"),
    @("801×7=5607", "607×6=3642"),
    @("404×8=3232", "179×2=358")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
